$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix mistake in pixel measurement: B7 was 67, should be 55
$ws.Range("B7").Value = 55

# Update active cell selection to C14
$ws.Range("C14").Select()
